$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 246.42857
$ws.Cells.Item(5, 9).Value = 314.4
$ws.Cells.Item(5, 10).Value = 76.5
$ws.Cells.Item(5, 11).Value = 314.4
$ws.Cells.Item(5, 12).Value = 76.5
$ws.Cells.Item(5, 13).Value = -199.4
$ws.Cells.Item(5, 14).Value = -306.5

$ws.Cells.Item(33, 8).Value = 89.28570999999999
$ws.Cells.Item(33, 9).Value = 89.28570999999999
$ws.Cells.Item(33, 11).Value = 89.28570999999999
$ws.Cells.Item(33, 13).Value = 139.71429

$ws.Cells.Item(34, 8).Value = 2241.8572
$ws.Cells.Item(34, 9).Value = 2241.8572
$ws.Cells.Item(34, 11).Value = 2241.8572
$ws.Cells.Item(34, 13).Value = -2038.8572

$ws.Cells.Item(36, 8).Value = 2241.8572
$ws.Cells.Item(36, 9).Value = 2241.8572
$ws.Cells.Item(36, 11).Value = 2241.8572
$ws.Cells.Item(36, 13).Value = -1526.8572

$ws.Cells.Item(70, 8).Value = 3070
$ws.Cells.Item(70, 9).Value = 3337.5
$ws.Cells.Item(70, 11).Value = 10012.5
$ws.Cells.Item(70, 13).Value = -9742.5

$ws.Cells.Item(73, 8).Value = 3070
$ws.Cells.Item(73, 9).Value = 3337.5
$ws.Cells.Item(73, 11).Value = 10012.5
$ws.Cells.Item(73, 13).Value = -9076.5

$ws.Cells.Item(93, 8).Value = 57500
$ws.Cells.Item(93, 10).Value = 57500
$ws.Cells.Item(93, 12).Value = 57500
$ws.Cells.Item(93, 14).Value = -62492

$ws.Cells.Item(100, 8).Value = 2324.75
$ws.Cells.Item(100, 9).Value = 2324.75
$ws.Cells.Item(100, 11).Value = 2324.75
$ws.Cells.Item(100, 13).Value = -1783.75

$ws.Cells.Item(125, 8).Value = 747.5
$ws.Cells.Item(125, 10).Value = 747.5
$ws.Cells.Item(125, 12).Value = 6727.5
$ws.Cells.Item(125, 14).Value = -11647.5

$ws.Cells.Item(132, 8).Value = 2012.7
$ws.Cells.Item(132, 9).Value = 718.4286
$ws.Cells.Item(132, 10).Value = 5032.6665
$ws.Cells.Item(132, 11).Value = 2155.2858
$ws.Cells.Item(132, 12).Value = 15097.9995
$ws.Cells.Item(132, 13).Value = 374.7142000000003
$ws.Cells.Item(132, 14).Value = -20157.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(52, 8).Value = 0
$ws.Cells.Item(52, 9).Value = 0
$ws.Cells.Item(52, 11).Value = 0
$ws.Cells.Item(52, 13).Value = $null

$ws.Cells.Item(61, 8).Value = 5000
$ws.Cells.Item(61, 9).Value = 5000
$ws.Cells.Item(61, 11).Value = 5000
$ws.Cells.Item(61, 13).Value = -4788

$ws.Cells.Item(74, 8).Value = 3464.95
$ws.Cells.Item(74, 9).Value = 3464.95
$ws.Cells.Item(74, 11).Value = 3464.95
$ws.Cells.Item(74, 13).Value = -2590.95

$ws.Cells.Item(77, 8).Value = 3464.95
$ws.Cells.Item(77, 9).Value = 3464.95
$ws.Cells.Item(77, 11).Value = 17324.75
$ws.Cells.Item(77, 13).Value = -12956.75

$ws.Cells.Item(98, 8).Value = 17499.5
$ws.Cells.Item(98, 10).Value = 17499.5
$ws.Cells.Item(98, 12).Value = 17499.5
$ws.Cells.Item(98, 14).Value = -23489.5

$ws.Cells.Item(136, 8).Value = 5000
$ws.Cells.Item(136, 9).Value = 5000
$ws.Cells.Item(136, 11).Value = 15000
$ws.Cells.Item(136, 13).Value = -12450

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(19, 8).Value = 26500
$ws.Cells.Item(19, 10).Value = 45000
$ws.Cells.Item(19, 12).Value = 45000
$ws.Cells.Item(19, 14).Value = -45346

$ws.Cells.Item(20, 8).Value = 0
$ws.Cells.Item(20, 9).Value = 0
$ws.Cells.Item(20, 11).Value = 0
$ws.Cells.Item(20, 13).Value = $null

$ws.Cells.Item(33, 8).Value = 9673.666999999999
$ws.Cells.Item(33, 9).Value = 9673.666999999999
$ws.Cells.Item(33, 10).Value = 0
$ws.Cells.Item(33, 11).Value = 9673.666999999999
$ws.Cells.Item(33, 12).Value = 0
$ws.Cells.Item(33, 13).Value = $null
$ws.Cells.Item(33, 14).Value = -9337.666999999999

$ws.Cells.Item(134, 8).Value = 6472.364
$ws.Cells.Item(134, 9).Value = 5369.6
$ws.Cells.Item(134, 11).Value = 16108.8
$ws.Cells.Item(134, 13).Value = -13573.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(3, 8).Value = 8334.333000000001
$ws.Cells.Item(3, 9).Value = 5000
$ws.Cells.Item(3, 10).Value = 10001.5
$ws.Cells.Item(3, 11).Value = 5000
$ws.Cells.Item(3, 12).Value = 10001.5
$ws.Cells.Item(3, 13).Value = -4887
$ws.Cells.Item(3, 14).Value = -10227.5

$ws.Cells.Item(16, 8).Value = 939.5185
$ws.Cells.Item(16, 9).Value = 898.7308
$ws.Cells.Item(16, 11).Value = 898.7308
$ws.Cells.Item(16, 13).Value = -611.7308

$ws.Cells.Item(22, 8).Value = 775.7143
$ws.Cells.Item(22, 10).Value = 1000
$ws.Cells.Item(22, 12).Value = 1000
$ws.Cells.Item(22, 14).Value = -1700

$ws.Cells.Item(35, 8).Value = 5380.222
$ws.Cells.Item(35, 10).Value = 5274.5
$ws.Cells.Item(35, 12).Value = 5274.5
$ws.Cells.Item(35, 14).Value = -5862.5

$ws.Cells.Item(58, 8).Value = 558.0909
$ws.Cells.Item(58, 9).Value = 467.5
$ws.Cells.Item(58, 10).Value = 799.6667
$ws.Cells.Item(58, 11).Value = 467.5
$ws.Cells.Item(58, 12).Value = 799.6667
$ws.Cells.Item(58, 13).Value = -264.5
$ws.Cells.Item(58, 14).Value = -1205.6667

$ws.Cells.Item(74, 8).Value = 98500
$ws.Cells.Item(74, 10).Value = 98500
$ws.Cells.Item(74, 12).Value = 98500
$ws.Cells.Item(74, 14).Value = -100248

$ws.Cells.Item(77, 8).Value = 98500
$ws.Cells.Item(77, 10).Value = 98500
$ws.Cells.Item(77, 12).Value = 295500
$ws.Cells.Item(77, 14).Value = -304236

$ws.Cells.Item(86, 8).Value = 5499.5
$ws.Cells.Item(86, 9).Value = 5499.5
$ws.Cells.Item(86, 11).Value = 5499.5
$ws.Cells.Item(86, 13).Value = -4376.5

$ws.Cells.Item(88, 8).Value = 20535.75
$ws.Cells.Item(88, 10).Value = 20535.75
$ws.Cells.Item(88, 12).Value = 20535.75
$ws.Cells.Item(88, 14).Value = -21347.75

$ws.Cells.Item(89, 8).Value = 5499.5
$ws.Cells.Item(89, 9).Value = 5499.5
$ws.Cells.Item(89, 11).Value = 27497.5
$ws.Cells.Item(89, 13).Value = -21881.5

$ws.Cells.Item(91, 8).Value = 20535.75
$ws.Cells.Item(91, 10).Value = 20535.75
$ws.Cells.Item(91, 12).Value = 20535.75
$ws.Cells.Item(91, 14).Value = -23343.75

$ws.Cells.Item(113, 8).Value = 939.5185
$ws.Cells.Item(113, 9).Value = 898.7308
$ws.Cells.Item(113, 11).Value = 898.7308
$ws.Cells.Item(113, 13).Value = 1271.2692

$ws.Cells.Item(134, 8).Value = 1369
$ws.Cells.Item(134, 9).Value = 1369
$ws.Cells.Item(134, 11).Value = 4107
$ws.Cells.Item(134, 13).Value = -1572

$ws.Cells.Item(135, 8).Value = 150000
$ws.Cells.Item(135, 10).Value = 150000
$ws.Cells.Item(135, 12).Value = 150000
$ws.Cells.Item(135, 14).Value = -160140

$ws.Cells.Item(136, 8).Value = 558.0909
$ws.Cells.Item(136, 9).Value = 467.5
$ws.Cells.Item(136, 10).Value = 799.6667
$ws.Cells.Item(136, 11).Value = 1402.5
$ws.Cells.Item(136, 12).Value = 2399.0001
$ws.Cells.Item(136, 13).Value = 1147.5
$ws.Cells.Item(136, 14).Value = -7499.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 189.08333
$ws.Cells.Item(2, 9).Value = 21.166666
$ws.Cells.Item(2, 10).Value = 357
$ws.Cells.Item(2, 11).Value = 126.999996
$ws.Cells.Item(2, 12).Value = 2142
$ws.Cells.Item(2, 13).Value = -13.999996
$ws.Cells.Item(2, 14).Value = -2368

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 5004.3335
$ws.Cells.Item(122, 9).Value = 5004.3335
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 15013.0005
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = $null
$ws.Cells.Item(122, 14).Value = -12563.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2005
$ws.Cells.Item(7, 10).Value = 2005
$ws.Cells.Item(7, 12).Value = 2005
$ws.Cells.Item(7, 14).Value = -2229

$ws.Cells.Item(40, 8).Value = 50004.5
$ws.Cells.Item(40, 9).Value = 10004
$ws.Cells.Item(40, 11).Value = 10004
$ws.Cells.Item(40, 13).Value = -9868

$ws.Cells.Item(46, 8).Value = 287785.16
$ws.Cells.Item(46, 10).Value = 2416
$ws.Cells.Item(46, 12).Value = 2416
$ws.Cells.Item(46, 14).Value = -2792

$ws.Cells.Item(104, 8).Value = 27249
$ws.Cells.Item(104, 10).Value = 27249
$ws.Cells.Item(104, 12).Value = 27249
$ws.Cells.Item(104, 14).Value = -34237

$ws.Cells.Item(122, 8).Value = 6749.75
$ws.Cells.Item(122, 9).Value = 5999.5
$ws.Cells.Item(122, 10).Value = 7500
$ws.Cells.Item(122, 11).Value = 17998.5
$ws.Cells.Item(122, 12).Value = 22500
$ws.Cells.Item(122, 13).Value = -15548.5
$ws.Cells.Item(122, 14).Value = -27400

$ws.Cells.Item(126, 8).Value = 2005
$ws.Cells.Item(126, 10).Value = 2005
$ws.Cells.Item(126, 12).Value = 6015
$ws.Cells.Item(126, 14).Value = -10955

$ws.Cells.Item(132, 8).Value = 6000
$ws.Cells.Item(132, 9).Value = 6000
$ws.Cells.Item(132, 11).Value = 18000
$ws.Cells.Item(132, 13).Value = -15470

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(5, 8).Value = 1858428.6
$ws.Cells.Item(5, 10).Value = 1858428.6
$ws.Cells.Item(5, 12).Value = 1858428.6
$ws.Cells.Item(5, 14).Value = -1858652.6

$ws.Cells.Item(17, 8).Value = 1661
$ws.Cells.Item(17, 10).Value = 2135
$ws.Cells.Item(17, 12).Value = 2135
$ws.Cells.Item(17, 14).Value = -2479

$ws.Cells.Item(80, 8).Value = 22666.334
$ws.Cells.Item(80, 10).Value = 23999.5
$ws.Cells.Item(80, 12).Value = 23999.5
$ws.Cells.Item(80, 14).Value = -25995.5

$ws.Cells.Item(83, 8).Value = 22666.334
$ws.Cells.Item(83, 10).Value = 23999.5
$ws.Cells.Item(83, 12).Value = 71998.5
$ws.Cells.Item(83, 14).Value = -81982.5
